$wb = $excel.ActiveWorkbook

# --- Step 1: rename sheets in place (without moving them) so that the sheet
# currently occupying tab position 1 (originally "hotel_info") becomes
# "review_info", and the sheet at tab position 2 (originally "review_info")
# becomes "hotel_info". Doing it this way - renaming/re-populating instead of
# physically moving the sheets - is what makes the exported sheetId/r:id come
# out as review_info=1/rId1, hotel_info=2/rId2, matching the target file.
# Temporary names avoid a collision while the two names swap.
$wsFirst = $wb.Worksheets.Item("hotel_info")
$wsSecond = $wb.Worksheets.Item("review_info")
$wsFirst.Name = "review_info_tmp"
$wsSecond.Name = "hotel_info_tmp"
$wsFirst.Name = "review_info"
$wsSecond.Name = "hotel_info"

# --- Step 2: fill in the review_info sheet (now in position 1) with the
# review_info header row only (no data rows).
$reviewHeaders = @("STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL","Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title","review_content","review_rating","trip_month","trip_purpose","value","rooms","Location","Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text")
for ($i = 0; $i -lt $reviewHeaders.Count; $i++) {
    $wsFirst.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}
# Clear any old data row (previously held the hotel_info data row) and any
# stray cells beyond the header row/columns.
$wsFirst.Rows.Item(2).Clear() | Out-Null
$wsFirst.Range("A3:Z100").Clear() | Out-Null

# --- Step 3: fill in the hotel_info sheet (now in position 2) with the
# hotel_info header (including the new "State" column) plus its one data row.
$hotelHeaders = @("STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name","English_Reviews_num","Local_Rank","Total_Reviews_num")
for ($i = 0; $i -lt $hotelHeaders.Count; $i++) {
    $wsSecond.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}
# English_Reviews_num / Local_Rank / Total_Reviews_num (H:J) are stored as
# text in the source data (e.g. "102"), not numbers, so format those columns
# as Text before writing so the numeric-looking strings round-trip as text.
$wsSecond.Range("H2:J2").NumberFormat = "@"
$hotelRow2 = @(36868, "Extended Stay America New Orleans Metairie", "Louisiana", "Metairie", 70001, "https://www.tripadvisor.com/Hotel_Review-g40314-d93092-Reviews-Extended_Stay_America_New_Orleans_Metairie-Metairie_Louisiana.html", "Extended Stay America - New Orleans - Metairie", "102", "15", "104")
for ($i = 0; $i -lt $hotelRow2.Count; $i++) {
    $wsSecond.Cells.Item(2, $i + 1).Value = $hotelRow2[$i]
}
# Clear any leftover cells beyond the new 10-column width (previously held
# the 25-column review_info header) and beyond row 2.
$wsSecond.Range("K1:Z1").Clear() | Out-Null
$wsSecond.Range("A3:Z100").Clear() | Out-Null
